$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '95.286.87'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.95%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.578.83'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.28%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.66'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.27%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '651.72'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.34%  '
$ws.Range("E8").Value = '  -0.42%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.00'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.00%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.578.77'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.22%  '
$ws.Range("E12").Value = '  +1.05%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '42.44'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.40%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.46'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.53%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.242.91'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.32%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '95.157.18'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.99%  '
$ws.Range("E17").Value = '  -0.35%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.561.69'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.78%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.75'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.36%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.54'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.21%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.85'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.70%  '
$ws.Range("E22").Value = '  +0.44%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '508.74'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.45%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.479'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.47%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.84'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.63%  '
$ws.Range("E26").Value = '  -2.27%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '95.34'
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.70'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.12%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.771.66'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.20%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.01'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.28%  '
$ws.Range("E31").Value = '  -1.00%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.47'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.55%  '
$ws.Range("E33").Value = '  +0.35%  '
$ws.Range("E34").Value = '  +0.90%  '
$ws.Range("E35").Value = '  -2.34%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '31.80'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.36%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.68'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +13.13%  '
$ws.Range("E38").Value = '  -1.03%  '
$ws.Range("E39").Value = '  +8.21%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '585.30'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.17%  '
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.151'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.86%  '
$ws.Range("E43").Value = '  -2.51%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.80'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.76%  '
$ws.Range("B45").Value = 'Stacks'
$ws.Range("C45").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.31'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +6.03%  '
$ws.Range("B46").Value = 'Filecoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.75'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.09%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '33.89'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +29.39%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.39'
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0416'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.72%  '
$ws.Range("E50").Value = '  -0.53%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '53.20'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.05%  '
